# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Row 4: Estados Unidos (United States) - updated totals ---
$ws.Cells.Item(4, 2).Value = 1363126   # B4 Casos totales
$ws.Cells.Item(4, 3).Value = 15817     # C4 Nuevos casos
$ws.Cells.Item(4, 5).Value = 1041964   # E4 Recuperados
$ws.Cells.Item(4, 7).Value = 509       # G4 Casos criticos
$ws.Cells.Item(4, 8).Value = 80546     # H4 Muertes

# --- Row 91: Republica de Yibuti - updated totals ---
$ws.Cells.Item(91, 2).Value = 1210     # B91 Casos totales
$ws.Cells.Item(91, 3).Value = 21       # C91 Nuevos casos
$ws.Cells.Item(91, 4).Value = 847      # D91 Casos activos
$ws.Cells.Item(91, 5).Value = 360      # E91 Recuperados

# --- Row 109: Principado de Andorra - updated totals ---
$ws.Cells.Item(109, 2).Value = 755     # B109 Casos totales
$ws.Cells.Item(109, 3).Value = 1       # C109 Nuevos casos
$ws.Cells.Item(109, 4).Value = 550     # D109 Casos activos
$ws.Cells.Item(109, 5).Value = 157     # E109 Recuperados

# --- Rows 192/193: Belice & Nueva Caledonia swap order/data ---
# Row 192 now shows Nueva Caledonia's data, row 193 now shows Belice's data
$ws.Cells.Item(192, 1).Value = "Nueva Caledonia"
$ws.Cells.Item(192, 4).Value = 18      # D192 Casos activos
$ws.Cells.Item(192, 8).Value = 0       # H192 Muertes

$ws.Cells.Item(193, 1).Value = "Belice"
$ws.Cells.Item(193, 4).Value = 16      # D193 Casos activos
$ws.Cells.Item(193, 8).Value = 2       # H193 Muertes

# --- Rows 212/213: Butan & Islas Virgenes Britanicas swap order/data ---
# Row 212 now shows Islas Virgenes Britanicas's data, row 213 now shows Butan's data
$ws.Cells.Item(212, 1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(212, 4).Value = 4       # D212 Casos activos
$ws.Cells.Item(212, 8).Value = 1       # H212 Muertes

$ws.Cells.Item(213, 1).Value = "Butan"
$ws.Cells.Item(213, 4).Value = 5       # D213 Casos activos
$ws.Cells.Item(213, 8).Value = 0       # H213 Muertes
